$wb = $excel.ActiveWorkbook

# Delete the now-unused "placeNames" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("placeNames").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheet to match the new naming convention
$ws = $wb.Worksheets.Item("yearNames")
$ws.Name = "cdli_years"

# Scroll back to the top and select cell B6
$ws.Range("B6").Select()
